$wb = $excel.ActiveWorkbook

# --- 1. Clear the stray empty cells B3 and B4 on "ODI Batting" sheet ---
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("B3").ClearContents()
$odiBatting.Range("B4").ClearContents()

# --- 2. Add a new worksheet "ODI Batting Extra" at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add($null, $lastSheet)
$extra.Name = "ODI Batting Extra"

# All cells on this sheet are text-typed in the source data (except the
# numeric BATTING_POSITION value in B3), so force the whole used range to
# text format before writing any values.
$extra.Range("A1:F5").NumberFormat = "@"

# Header row
$extra.Range("A1").Value = "MATCH_CODE"
$extra.Range("B1").Value = "BATTING_POSITION"
$extra.Range("C1").Value = "NUM_4"
$extra.Range("D1").Value = "NUM_6"
$extra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Range("F1").Value = "MAN_OF_MATCH"

# Row 2 - match 4564
$extra.Range("A2").Value = "4564"
$extra.Range("F2").Value = "NO"

# Row 3 - match 4565
$extra.Range("A3").Value = "4565"
$extra.Range("B3").NumberFormat = "General"
$extra.Range("B3").Value = 10
$extra.Range("C3").Value = "0"
$extra.Range("D3").Value = "1"
$extra.Range("E3").Value = "4.37%"
$extra.Range("F3").Value = "NO"

# Row 4 - match 4567 (only MATCH_CODE is populated)
$extra.Range("A4").Value = "4567"

# Row 5 - match 4641 (only MATCH_CODE is populated)
$extra.Range("A5").Value = "4641"

# Header style: bold, centered, bordered (same style used by the other sheets)
$extra.Range("A1:F1").Font.Bold = $true
$extra.Range("A1:F1").HorizontalAlignment = -4108
$extra.Range("A1:F1").VerticalAlignment = -4160
$extra.Range("A1:F1").Borders.LineStyle = 1
